# OW-268 - fix data format in the excel files to follow strata conventions.
#
# The "IRS-Cleared" sheet carries two convention-named fields for both legs
# of the swap:
#   - LEG1_PAY_ADJ_BUS_DAY_CONV / LEG2_PAY_ADJ_BUS_DAY_CONV (cols S / AF)
#       "MODFOLLOWING"    -> "ModifiedFollowing"
#   - LEG1_DAYCOUNT       / LEG2_DAYCOUNT       (cols U / AH)
#       "ACT/365.FIXED"   -> "Act/365F"
#
# Update the four data cells on row 2 to the Strata-style convention names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Value = "ModifiedFollowing"
$ws.Range("U2").Value = "Act/365F"
$ws.Range("AF2").Value = "ModifiedFollowing"
$ws.Range("AH2").Value = "Act/365F"

# Leave the reviewer's cursor on the day-count cell that was touched last,
# mirroring where focus ended up after making the edit.
[void]$ws.Range("U2").Select()
